# KIBON-2685 attribute zu excel export hinzufuegen
# Extends the LastenausgleichBGZeitabschnitte.xlsx export template with the
# per-Betreuung detail columns (Referenznummer, BFS-Nummer, Gemeinde, Nachname,
# Vorname, Geburtsdatum, von/bis, Institution, Betreuungsangebot, BG-Pensum,
# kein Selbstbehalt, Gutschein) plus a title block (Lastenausgleich-Daten /
# Parameter / Jahr) above the existing repeating-row table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1: big bold title
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "{lastenausgleichDatenTitle}"
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").Font.Size = 16
$ws.Rows("1:1").RowHeight = 21

# ---------------------------------------------------------------------------
# Row 3/4: parameter / jahr block
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "{parameterTitle}"
$ws.Range("A3:B3").Font.Bold = $true

$ws.Range("A4").Value = "{jahrTitle}"
$ws.Range("B4").Value = "{jahr}"
$ws.Range("B4").NumberFormat = "0"
$ws.Range("B4").Interior.ThemeColor = 3
$ws.Range("B4").Interior.TintAndShade = 0.8
$ws.Range("B4").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# Row 7 (header, bold/wrapped/grey) and row 8 (placeholder values)
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "{referenznummerTitle}"
$ws.Range("A8").Value = "{referenznummer}"

$ws.Range("B7").Value = "{bfsNummerTitle}"
$ws.Range("B8").Value = "{bfsNummer}"

$ws.Range("C7").Value = "{nameGemeindeTitle}"
$ws.Range("C8").Value = "{nameGemeinde}"

$ws.Range("D7").Value = "{nachnameTitle}"
$ws.Range("D8").Value = "{nachname}"

$ws.Range("E7").Value = "{vornameTitle}"
$ws.Range("E8").Value = "{vorname}"

$ws.Range("F7").Value = "{geburtsdatumTitle}"
$ws.Range("F8").Value = "{geburtsdatum}"
$ws.Range("F7:F8").NumberFormat = "mm-dd-yy"

$ws.Range("G7").Value = "{vonTitle}"
$ws.Range("G8").Value = "{von}"
$ws.Range("G7:G8").NumberFormat = "mm-dd-yy"

$ws.Range("H7").Value = "{bisTitle}"
$ws.Range("H8").Value = "{bis}"
$ws.Range("H7:H8").NumberFormat = "mm-dd-yy"

$ws.Range("I7").Value = "{institutionTitle}"
$ws.Range("I8").Value = "{institution}"

$ws.Range("J7").Value = "{betreuungsangebotTypTitle}"
$ws.Range("J8").Value = "{betreuungsangebotTyp}"

$ws.Range("K7").Value = "{bgPensumTitle}"
$ws.Range("K8").Value = "{bgPensum}"
$ws.Range("K7:K8").NumberFormat = "0.00%"

$ws.Range("L7").Value = "{keinSelbstbehaltDurchGemeindeTitle}"
$ws.Range("L8").Value = "{keinSelbstbehaltDurchGemeinde}"

$ws.Range("M7").Value = "{gutscheinTitle}"
$ws.Range("M8").Value = "{gutschein}"
$ws.Range("M7:M8").NumberFormat = """CHF"" #,##0.00"

$ws.Range("N8").Value = "{repeatRow}"

# Header row formatting (bold/wrap/grey fill/border), matches existing A1(old)/A7 style
$ws.Range("A7:M7").Font.Bold = $true
$ws.Range("A7:M7").WrapText = $true
$ws.Range("A7:M7").Interior.Pattern = -4124
$ws.Range("A7:M7").Interior.ThemeColor = 1
$ws.Range("A7:M7").Interior.TintAndShade = -0.0499893185216834
$ws.Range("A7:M7").Borders.LineStyle = 1
$ws.Rows("7:7").RowHeight = 57.6

$ws.Range("A8:M8").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# Column widths (bestFit columns added for the new data columns)
# ---------------------------------------------------------------------------
$ws.Range("B:B").ColumnWidth = 11
$ws.Range("C:C").ColumnWidth = 11.88671875
$ws.Range("D:D").ColumnWidth = 18.109375
$ws.Range("I:I").ColumnWidth = 8.109375
$ws.Range("J:J").ColumnWidth = 11.33203125
$ws.Range("L:L").ColumnWidth = 11.33203125
$ws.Range("N:N").ColumnWidth = 11.44140625
